# Sun Dec 25 15:32:53 UTC 2022 symbol-list refresh (GitHub Actions bot).
# The source sheet is an auto-generated coin ranking snapshot where every
# cell -- including the "Price" column -- is stored as text (t="inlineStr"),
# so we replicate that by forcing a Text number format before writing any
# value that looks numeric (otherwise Excel would silently coerce it to a
# real number and trailing zeros like "244.20" / "1.050" would be lost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link / Volume(1h) columns: plain text, safe to assign directly ---
$textUpdates = @{
    'B7' = 'MXToken'
    'C7' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'E7' = '6MXTokenMX'
    'B8' = 'FTXToken'
    'C8' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'E8' = '7FTXTokenFTT'
    'B9' = 'WazirX'
    'C9' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'E9' = '8WazirXWRX'
    'B10' = 'MandalaExchangeToken'
    'C10' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'E10' = '9MandalaExchangeTokenMDX'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'E11' = '10LiechtensteinCryptoassetsExchangeLCX'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'E12' = '11BitrueCoinBTR'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'E13' = '12BitMartTokenBMX'
    'B14' = 'MCDex'
    'C14' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'E14' = '13MCDexMCB'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'E15' = '14BitForexTokenBF'
    'B16' = 'CoinExToken'
    'C16' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'E16' = '15CoinExTokenCET'
    'B17' = 'One'
    'C17' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E17' = '16OneONE'
    'B18' = 'TigerCash'
    'C18' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'E18' = '17TigerCashTCH'
    'B19' = 'HotbitToken'
    'C19' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'E19' = '18HotbitTokenHTB'
    'B20' = 'BitKan'
    'C20' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'E20' = '19BitKanKAN'
    'B21' = 'NitroEx'
    'C21' = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
    'E21' = '20NitroExNTX'
    'B22' = 'LEO'
    'C22' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'E22' = '21LEOLEO'
    'B23' = 'KuCoinToken'
    'C23' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'E23' = '22KuCoinTokenKCS'
    'E41' = '40KickTokenKICK'
    'E48' = '47CoinbaseStockTokenCOINBestin24h'
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}

# --- Price column: numeric-looking strings that must stay text ---
$priceUpdates = @{
    'D2' = '244.20'
    'D4' = '5.414'
    'D5' = '0.05932'
    'D6' = '3.402'
    'D7' = '0.8071'
    'D8' = '0.9162'
    'D9' = '0.1421'
    'D10' = '0.07424'
    'D11' = '0.03328'
    'D12' = '0.03076'
    'D13' = '0.09344'
    'D14' = '3.951'
    'D15' = '0.001587'
    'D16' = '0.04772'
    'D17' = '0.0005944'
    'D18' = '0.005470'
    'D19' = '0.004441'
    'D20' = '0.0009860'
    'D21' = '0.00007804'
    'D22' = '3.658'
    'D23' = '6.437'
    'D26' = '0.1340'
    'D40' = '0.03897'
    'D41' = '0.006216'
    'D42' = '0.1071'
    'D43' = '0.002611'
    'D44' = '0.006511'
    'D45' = '0.00005202'
    'D46' = '0.00000000751'
    'D47' = '0.0005804'
    'D48' = '1.050'
    'D49' = '0.002274'
    'D50' = '0.00002102'
    'D51' = '0.0002001'
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"            # force text storage
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = 'Normal'              # drop the temporary text format again
}

